# Applies the "ACTUALLY FIXED" commit:
#  - Drop the trailing "AI Insights" / "AI Recommendations" / "Statistical
#    Summary" slides (3, 4, 5) - they're reverted back out of the deck,
#    along with their bound notes pages.
#  - Slide 1 title: "test-report.xlsx" -> "test.xlsx"
#  - Slide 2 body: "This is a test summary for the AI report analysis."
#    -> "Test summary"

$p = $ppt.ActivePresentation

# --- Drop slides 3, 4, 5 (highest index first so indices stay valid) ---
for ($i = $p.Slides.Count; $i -ge 3; $i--) {
    $p.Slides.Item($i).Delete()
}

# --- Slide 1: "test-report.xlsx" -> "test.xlsx" ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item("Text 0").TextFrame.TextRange.Text = "test.xlsx"

# --- Slide 2: summary placeholder text -> "Test summary" ---
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item("Text 1").TextFrame.TextRange.Text = "Test summary"
